# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font + border + centered/top alignment) from
# the existing header row into the three new header cells, then set text.
$ws.Range("A1:C1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill the season record for every player row (2-56) with the team's
# overall record for the season: 69 wins, 93 losses, 0 ties.
$ws.Range("AD2:AD56").Value = 69
$ws.Range("AE2:AE56").Value = 93
$ws.Range("AF2:AF56").Value = 0
